# VerveStacks ITA model update - 2025-09-01 16:30
# Adds a new scenario-group ("s1_d" / "f3d") to the ScenMap lookup table and
# extends the generated scenario list (rows 56:60) to match it, then leaves
# the ScenMap sheet active with A11 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenMap")

# --- 1. Register the new lookup-table row (S16:U16) -----------------------
$ws.Range("S16").Value = 11
$ws.Range("T16").Value = "s1_d"
$ws.Range("U16").Value = "f3d"

# --- 2. Append the five generated rows (56:60) for the new group ----------
$ws.Range("A56:A60").Formula = "=A51+1"
$ws.Range("B56:B60").Formula = '="vstacks_"&VLOOKUP(A56,$S$6:$T$18,2,FALSE)&"~"&TEXT(O56,"0000")'
$ws.Range("C56:C60").Formula = "=H56"
$ws.Range("H56:H60").Formula = '=TEXTJOIN(".",TRUE,I56:J56)'
$ws.Range("I56:I60").Formula = "=I51"
$ws.Range("J56:J60").Formula = "=Q56"
$ws.Range("O56:O60").Formula = "=O51"
$ws.Range("Q56:Q60").Formula = '=VLOOKUP(A56,$S$6:$U$17,3,FALSE)'

# --- 3. Make ScenMap the active sheet/selection, matching the saved view --
$ws.Activate()
$ws.Range("A11").Select()
